$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price strings to remain plain text (matching the
# source feed formatting) instead of being auto-converted to Number by Excel.
$textCells = @("D5","D6","D12","D13","D14","D19","D20","D21","D22","D23","D24","D25","D27","D30","D31","D32","D37","D38","D39","D40","D41","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated cryptos list values per the scheduled data refresh.
# Rows 39 and 40 swap (Stacks <-> EthereumClassic) in addition to value updates.
$ws.Range("D2").Value = "68.814.93"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.651.60"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "601.00"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "155.78"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "2.649.12"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +14.09%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "27.96"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  +6.39%  "
$ws.Range("D16").Value = "3.134.43"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "68.736.25"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "2.644.01"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "11.44"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").Value = "367.74"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "7.48"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("D25").Value = "73.08"
$ws.Range("E25").Value = "  +10.50%  "
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  +6.43%  "
$ws.Range("D29").Value = "2.779.65"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "583.03"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "1.43"
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("E33").Value = "  +4.59%  "
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "1.56"
$ws.Range("E37").Value = "  +4.21%  "
$ws.Range("D38").Value = "160.11"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "19.36"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("D41").Value = "5.42"
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("E43").Value = "  +5.68%  "
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("D45").Value = "0.0₆0320"
$ws.Range("E45").Value = "  +9.28%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "40.59"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "156.73"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").Value = "3.75"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "22.12"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  +1.27%  "
